# Applies the crypto price/volume/ranking refresh described by the
# commit "Updated cryptos list on Mon May 22 17:41:56 UTC 2023 with GitHub Actions".
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). All four columns hold text
# (not numbers), so every write below is forced to stay text: values that
# look numeric (e.g. "312.51") are written with a leading apostrophe so
# Excel doesn't coerce them into a Double, and ClearFormats() immediately
# strips the transient quotePrefix style flag that the apostrophe trick
# leaves behind, so the cell's style index is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $range = $ws.Range($cellAddr)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
    $range.ClearFormats()
}


# Row 2
Set-TextValue "D2" '26.988.26'
Set-TextValue "E2" '  -0.52%  '

# Row 3
Set-TextValue "D3" '1.826.84'
Set-TextValue "E3" '  +0.23%  '

# Row 4
Set-TextValue "E4" '  -0.37%  '

# Row 5
Set-TextValue "D5" '312.51'
Set-TextValue "E5" '  +0.04%  '

# Row 6
Set-TextValue "D6" '1.005'
Set-TextValue "E6" '  -0.48%  '

# Row 7
Set-TextValue "D7" '0.4569'
Set-TextValue "E7" '  -0.94%  '

# Row 8
Set-TextValue "D8" '0.3701'
Set-TextValue "E8" '  +1.88%  '

# Row 9
Set-TextValue "D9" '0.07353'
Set-TextValue "E9" '  +0.78%  '

# Row 10
Set-TextValue "D10" '0.8767'
Set-TextValue "E10" '  +0.65%  '

# Row 11
Set-TextValue "B11" 'TRON'
Set-TextValue "C11" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D11" '0.07945'
Set-TextValue "E11" '  +4.04%  '

# Row 12
Set-TextValue "B12" 'Solana'
Set-TextValue "C12" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue "D12" '19.76'
Set-TextValue "E12" '  -1.56%  '

# Row 13
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.756.28'
Set-TextValue "E13" '  -6.57%  '

# Row 14
Set-TextValue "B14" 'Chainlink'
Set-TextValue "C14" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D14" '6.602'
Set-TextValue "E14" '  +1.79%  '

# Row 15
Set-TextValue "B15" 'Polkadot'
Set-TextValue "C15" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D15" '5.337'
Set-TextValue "E15" '  -0.04%  '

# Row 16
Set-TextValue "D16" '91.49'
Set-TextValue "E16" '  -0.96%  '

# Row 17
Set-TextValue "D17" '1.008'
Set-TextValue "E17" '  -0.22%  '

# Row 18
Set-TextValue "D18" '0.000008934'
Set-TextValue "E18" '  +3.40%  '

# Row 19
Set-TextValue "D19" '1.006'
Set-TextValue "E19" '  -0.47%  '

# Row 20
Set-TextValue "B20" 'WrappedBTC'
Set-TextValue "C20" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D20" '27.728.88'
Set-TextValue "E20" '  +1.11%  '

# Row 21
Set-TextValue "B21" 'Avalanche'
Set-TextValue "C21" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D21" '14.82'
Set-TextValue "E21" '  +2.29%  '

# Row 22
Set-TextValue "D22" '5.113'
Set-TextValue "E22" '  -1.84%  '

# Row 23
Set-TextValue "B23" 'WrappedliquidstakedEther2.0'
Set-TextValue "C23" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D23" '2.246.24'
Set-TextValue "E23" '  +7.13%  '

# Row 24
Set-TextValue "B24" 'Cosmos'
Set-TextValue "C24" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D24" '10.53'
Set-TextValue "E24" '  -0.33%  '

# Row 25
Set-TextValue "D25" '153.17'
Set-TextValue "E25" '  +1.12%  '

# Row 26
Set-TextValue "D26" '1.847'
Set-TextValue "E26" '  -1.27%  '

# Row 27
Set-TextValue "D27" '18.37'
Set-TextValue "E27" '  +0.85%  '

# Row 28
Set-TextValue "D28" '2.044'
Set-TextValue "E28" '  -1.41%  '

# Row 29
Set-TextValue "D29" '5.143'
Set-TextValue "E29" '  +0.80%  '

# Row 30
Set-TextValue "D30" '115.32'
Set-TextValue "E30" '  -0.70%  '

# Row 31
Set-TextValue "D31" '0.08872'
Set-TextValue "E31" '  -0.39%  '

# Row 32
Set-TextValue "D32" '2.964'
Set-TextValue "E32" '  -0.02%  '

# Row 33
Set-TextValue "D33" '0.7295'
Set-TextValue "E33" '  -0.73%  '

# Row 34
Set-TextValue "D34" '4.419'
Set-TextValue "E34" '  -0.80%  '

# Row 35
Set-TextValue "D35" '1.131'
Set-TextValue "E35" '  -0.48%  '

# Row 36
Set-TextValue "B36" 'RenderToken'
Set-TextValue "C36" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D36" '2.460'
Set-TextValue "E36" '  -0.67%  '

# Row 37
Set-TextValue "B37" 'TrustWalletToken'
Set-TextValue "C37" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D37" '1.072'
Set-TextValue "E37" '  -0.19%  '

# Row 38
Set-TextValue "D38" '0.01943'
Set-TextValue "E38" '  +1.48%  '

# Row 39
Set-TextValue "D39" '0.05217'
Set-TextValue "E39" '  -0.59%  '

# Row 40
Set-TextValue "D40" '2.934'
Set-TextValue "E40" '  +0.40%  '

# Row 41
Set-TextValue "D41" '7.172'
Set-TextValue "E41" '  +0.20%  '

# Row 42
Set-TextValue "B42" 'TheSandbox'
Set-TextValue "C42" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D42" '0.5140'
Set-TextValue "E42" '  -1.02%  '

# Row 43
Set-TextValue "B43" 'Frax'
Set-TextValue "C43" 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue "D43" '0.8906'
Set-TextValue "E43" '  -12.01%  '

# Row 44
Set-TextValue "D44" '0.1628'
Set-TextValue "E44" '  +0.14%  '

# Row 45
Set-TextValue "D45" '8.180'
Set-TextValue "E45" '  -1.20%  '

# Row 46
Set-TextValue "D46" '0.4830'
Set-TextValue "E46" '  -0.18%  '

# Row 47
Set-TextValue "D47" '1.006'
Set-TextValue "E47" '  -0.54%  '

# Row 48
Set-TextValue "D48" '10.19'
Set-TextValue "E48" '  +0.11%  '

# Row 49
Set-TextValue "D49" '102.75'
Set-TextValue "E49" '  -0.50%  '

# Row 50
Set-TextValue "D50" '1.631'
Set-TextValue "E50" '  -0.14%  '

# Row 51
Set-TextValue "E51" '  -0.94%  '
